$wb = $excel.ActiveWorkbook

# Sheet "M" is the first sheet (sheet1.xml) - update period column (A2:A7)
# from plain years to monthly period labels.
$wsM = $wb.Worksheets.Item("M")
$wsM.Range("A2").Value = "1998M01"
$wsM.Range("A3").Value = "1998M02"
$wsM.Range("A4").Value = "1998M03"
$wsM.Range("A5").Value = "1998M04"
$wsM.Range("A6").Value = "1998M05"
$wsM.Range("A7").Value = "1998M06"

# Move the active tab / selection from sheet "A" to sheet "M",
# and update the selected range on sheet "M" to A8.
$wsM.Activate()
$wsM.Range("A8").Select()
